$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 24.999262
$ws.Range("H2").Value = 74.997786
$ws.Range("I2").Value = 0.2094245171924971
$ws.Range("J2").Value = 0.209424517192497
$ws.Range("M2").Value = 5.353120999999999
$ws.Range("N2").Value = 16.059363
$ws.Range("O2").Value = 0.3576717770696169
$ws.Range("P2").Value = 0.3576717770696169
$ws.Range("Q2").Value = 133.824074396702
$ws.Range("R2").Value = 1204.416669570318
$ws.Range("S2").Value = 0.07490523922618697
$ws.Range("T2").Value = 0.07490523922618696

$ws.Range("G3").Value = 24.999262
$ws.Range("H3").Value = 74.997786
$ws.Range("I3").Value = 0.2094245171924971
$ws.Range("J3").Value = 0.209424517192497
$ws.Range("O3").Value = 0.4477725819276249
$ws.Range("P3").Value = 0.4477725819276249
$ws.Range("Q3").Value = 167.5355875367893
$ws.Range("R3").Value = 1507.820287831104
$ws.Range("S3").Value = 0.09377455678223068
$ws.Range("T3").Value = 0.09377455678223065

$ws.Range("G4").Value = 24.999262
$ws.Range("H4").Value = 74.997786
$ws.Range("I4").Value = 0.2094245171924971
$ws.Range("J4").Value = 0.209424517192497
$ws.Range("O4").Value = 0.1945556410027581
$ws.Range("P4").Value = 0.1945556410027581
$ws.Range("Q4").Value = 72.79363440180934
$ws.Range("R4").Value = 655.142709616284
$ws.Range("S4").Value = 0.0407447211840794
$ws.Range("T4").Value = 0.0407447211840794

$ws.Range("I5").Value = 0.522807373179233
$ws.Range("J5").Value = 0.5228073731792329
$ws.Range("M5").Value = 5.353120999999999
$ws.Range("N5").Value = 16.059363
$ws.Range("O5").Value = 0.3576717770696169
$ws.Range("P5").Value = 0.3576717770696169
$ws.Range("Q5").Value = 334.0784247298653
$ws.Range("R5").Value = 3006.705822568788
$ws.Range("S5").Value = 0.1869934422301147
$ws.Range("T5").Value = 0.1869934422301146

$ws.Range("I6").Value = 0.522807373179233
$ws.Range("J6").Value = 0.5228073731792329
$ws.Range("O6").Value = 0.4477725819276249
$ws.Range("P6").Value = 0.4477725819276249
$ws.Range("S6").Value = 0.2340988073392645
$ws.Range("T6").Value = 0.2340988073392644

$ws.Range("I7").Value = 0.522807373179233
$ws.Range("J7").Value = 0.5228073731792329
$ws.Range("O7").Value = 0.1945556410027581
$ws.Range("P7").Value = 0.1945556410027581
$ws.Range("S7").Value = 0.1017151236098539
$ws.Range("T7").Value = 0.1017151236098538

$ws.Range("I8").Value = 0.2677681096282701
$ws.Range("J8").Value = 0.2677681096282701
$ws.Range("M8").Value = 5.353120999999999
$ws.Range("N8").Value = 16.059363
$ws.Range("O8").Value = 0.3576717770696169
$ws.Range("P8").Value = 0.3576717770696169
$ws.Range("Q8").Value = 171.1061336291416
$ws.Range("R8").Value = 1539.955202662275
$ws.Range("S8").Value = 0.09577309561331537
$ws.Range("T8").Value = 0.09577309561331537

$ws.Range("I9").Value = 0.2677681096282701
$ws.Range("J9").Value = 0.2677681096282701
$ws.Range("O9").Value = 0.4477725819276249
$ws.Range("P9").Value = 0.4477725819276249
$ws.Range("S9").Value = 0.1198992178061298
$ws.Range("T9").Value = 0.1198992178061298

$ws.Range("I10").Value = 0.2677681096282701
$ws.Range("J10").Value = 0.2677681096282701
$ws.Range("O10").Value = 0.1945556410027581
$ws.Range("P10").Value = 0.1945556410027581
$ws.Range("S10").Value = 0.05209579620882489
$ws.Range("T10").Value = 0.05209579620882489
